$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "总计" sheet: insert a new row for 2022-Q4 right under the header,
#    shifting the existing quarters down by one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 8
$summary.Cells.Item(2,4).Value = 1.25

# Row 2 picked up row 1's (header) formatting from the insert - restore it
# to match the plain-data-row formatting used by the other rows.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

# The row that used to be "2020-Q4" (now row 7) needs its running index
# fixed up to 5 (0-based position in the table).
$summary.Cells.Item(7,1).Value = 5

# ------------------------------------------------------------------
# 2. Add the new "2022-Q4" detail sheet right after "总计".
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"
$newSheet.Move($wb.Worksheets.Item(2))

$ws = $wb.Worksheets.Item(2)
$donor = $wb.Worksheets.Item(3)   # "2022-Q1" sheet - donates matching cell formats

$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# B, D, E, F, G hold text data in this workbook (fund codes with leading
# zeros, percentages kept with their original trailing zeros) - a leading
# "'" forces them to stay text instead of being parsed as numbers.
$rows = @(
    @("010106", "华夏核心科技6个月定期开放混合A", "6.19", "90.58", "7.41", "0.4587", 4),
    @("002148", "国寿安保稳惠灵活配置混合", "6.96", "77.40", "5.48", "0.3814", 3),
    @("013323", "国寿安保盛泽三年持有期混合A", "2.88", "74.45", "4.33", "0.1247", 6),
    @("168002", "国寿安保策略精选灵活配置混合（LOF）", "2.09", "75.43", "5.03", "0.1051", 5),
    @("010107", "华夏核心科技6个月定期开放混合C", "1.08", "90.58", "7.41", "0.0800", 4),
    @("009917", "格林泓利增强债券C", "2.82", "24.08", "2.43", "0.0685", 5),
    @("009916", "格林泓利增强债券A", "0.97", "24.08", "2.43", "0.0236", 5),
    @("013324", "国寿安保盛泽三年持有期混合C", "0.13", "74.45", "4.33", "0.0056", 6)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = ($r - 2)
    $ws.Cells.Item($r,2).Value = "'" + $row[0]
    $ws.Cells.Item($r,3).Value = $row[1]
    $ws.Cells.Item($r,4).Value = "'" + $row[2]
    $ws.Cells.Item($r,5).Value = "'" + $row[3]
    $ws.Cells.Item($r,6).Value = "'" + $row[4]
    $ws.Cells.Item($r,7).Value = "'" + $row[5]
    $ws.Cells.Item($r,8).Value = $row[6]
    $r = $r + 1
}

# Copy over the donor sheet's cell formatting (bold/bordered header row and
# index column, plain data cells) so the new sheet matches the others.
$donor.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$donor.Range("A2:A5").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)

$donor.Range("B2:H2").Copy()
$ws.Range("B2:H9").PasteSpecial(-4122)
